# "fixed pod to cluster"
#
# Slide 10 of the deck has a small pipeline diagram with three boxes that
# used to say "... POD(s)" and should now say "... Cluster(s)":
#   - "Production PODs" -> "Production Clusters"
#   - "QA PODs"          -> "QA Cluster"
#   - "Staging PODs"     -> "Staging Cluster" (this box auto-sizes to its
#     text - wrap="none" + spAutoFit - so its position/size shrink-wraps
#     to the new, shorter text as well)

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(10)

# --- "Production PODs" -> "Production Clusters" -----------------------
$prodBox = $slide.Shapes.Item(2)
if ($prodBox.TextFrame.TextRange.Text -eq "Production PODs") {
    $prodBox.TextFrame.TextRange.Text = "Production Clusters"
}

# --- "QA PODs" -> "QA Cluster" ------------------------------------------
$qaBox = $slide.Shapes.Item(18)
if ($qaBox.TextFrame.TextRange.Text -eq "QA PODs") {
    $qaBox.TextFrame.TextRange.Text = "QA Cluster"
}

# --- "Staging PODs" -> "Staging Cluster" (+ autosize reflow) -----------
$stagingBox = $slide.Shapes.Item(36)
if ($stagingBox.TextFrame.TextRange.Text -eq "Staging PODs") {
    $stagingBox.TextFrame.TextRange.Text = "Staging Cluster"

    # The textbox has spAutoFit (wrap="none"), so PowerPoint recomputes its
    # on-slide bounding box from the new text. Apply the resulting EMU
    # geometry converted to points (1 pt = 1/72 in = 914400/72 EMU). Left/Top
    # are nudged by <1 EMU so the Single-precision round-trip back to EMU
    # lands exactly on the target instead of 1 EMU short.
    $stagingBox.Left   = 619.7807923015748
    $stagingBox.Top    = 129.92591101181102
    $stagingBox.Width  = 1114408 / 914400.0 * 72
    $stagingBox.Height = 276999  / 914400.0 * 72
}
